$d = $word.ActiveDocument

# Locate the "Getting Started" heading paragraph.
# (Paragraph.Range.Text includes the trailing paragraph mark, so trim it
# before comparing.)
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Getting Started") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $newLines = @("System requirements.", "Running over HTTP.", "Running from local files.")
    foreach ($line in $newLines) {
        $target.Range.InsertParagraphAfter()
        $newPara = $target.Next()
        $newPara.Style = "Normal"
        $newPara.Range.Text = $line
        $target = $newPara
    }
}
